$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.403.95'
$ws.Range('E2').Value = '  +1.83%  '

$ws.Range('D3').Value = '1.852.65'
$ws.Range('E3').Value = '  +1.09%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = "'245.22"
$ws.Range('E5').Value = '  -0.08%  '

$ws.Range('D6').Value = "'0.6922"
$ws.Range('E6').Value = '  +0.41%  '

$ws.Range('D8').Value = "'0.07662"
$ws.Range('E8').Value = '  -0.60%  '

$ws.Range('D9').Value = "'0.3060"
$ws.Range('E9').Value = '  +0.39%  '

$ws.Range('D10').Value = "'23.49"

$ws.Range('D11').Value = "'0.07755"
$ws.Range('E11').Value = '  -0.65%  '

$ws.Range('D12').Value = "'5.143"
$ws.Range('E12').Value = '  +1.05%  '

$ws.Range('D13').Value = '1.853.98'
$ws.Range('E13').Value = '  +1.08%  '

$ws.Range('D14').Value = "'0.6931"
$ws.Range('E14').Value = '  +1.70%  '

$ws.Range('D15').Value = "'90.95"
$ws.Range('E15').Value = '  -0.05%  '

$ws.Range('D16').Value = "'6.303"
$ws.Range('E16').Value = '  -2.00%  '

$ws.Range('D17').Value = '29.432.24'
$ws.Range('E17').Value = '  +1.85%  '

$ws.Range('D18').Value = "'0.000008273"
$ws.Range('E18').Value = '  -0.46%  '

$ws.Range('D19').Value = '2.103.16'
$ws.Range('E19').Value = '  +1.00%  '

$ws.Range('D20').Value = "'236.26"
$ws.Range('E20').Value = '  -2.48%  '

$ws.Range('D21').Value = "'12.70"
$ws.Range('E21').Value = '  -0.08%  '

$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('E23').Value = '  +2.60%  '

$ws.Range('E24').Value = '  +0.12%  '

$ws.Range('E25').Value = '  -0.25%  '

$ws.Range('D26').Value = "'8.933"
$ws.Range('E26').Value = '  +1.62%  '

$ws.Range('D27').Value = "'159.93"
$ws.Range('E27').Value = '  +0.84%  '

$ws.Range('E28').Value = '  -0.08%  '

$ws.Range('D29').Value = "'1.529"
$ws.Range('E29').Value = '  -0.96%  '

$ws.Range('D30').Value = "'4.245"
$ws.Range('E30').Value = '  +0.62%  '

$ws.Range('D31').Value = "'4.134"
$ws.Range('E31').Value = '  -0.37%  '

$ws.Range('D32').Value = "'1.202"
$ws.Range('E32').Value = '  +0.88%  '

$ws.Range('D33').Value = "'0.05234"
$ws.Range('E33').Value = '  +2.64%  '

$ws.Range('D34').Value = "'0.7730"
$ws.Range('E34').Value = '  -0.80%  '

$ws.Range('E35').Value = '  +0.90%  '

$ws.Range('D36').Value = "'1.144"
$ws.Range('E36').Value = '  +0.28%  '

$ws.Range('D37').Value = "'2.692"
$ws.Range('E37').Value = '  +0.07%  '

$ws.Range('D38').Value = '1.328.94'
$ws.Range('E38').Value = '  +8.76%  '

$ws.Range('D39').Value = "'0.01864"
$ws.Range('E39').Value = '  +0.65%  '

$ws.Range('D40').Value = "'2.721"
$ws.Range('E40').Value = '  +1.05%  '

$ws.Range('D41').Value = "'0.9409"
$ws.Range('E41').Value = '  -1.44%  '

$ws.Range('D42').Value = "'105.85"
$ws.Range('E42').Value = '  -2.75%  '

$ws.Range('D43').Value = "'5.806"
$ws.Range('E43').Value = '  -0.31%  '

$ws.Range('E44').Value = '  +0.15%  '

$ws.Range('D45').Value = "'9.698"
$ws.Range('E45').Value = '  +0.83%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = "'0.00000000124"
$ws.Range('E46').Value = '  +1.39%  '

$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.999.28'
$ws.Range('E47').Value = '  +1.11%  '

$ws.Range('D48').Value = "'0.5226"
$ws.Range('E48').Value = '  +1.37%  '

$ws.Range('D49').Value = "'1.782"
$ws.Range('E49').Value = '  +1.93%  '

$ws.Range('D50').Value = "'62.90"
$ws.Range('E50').Value = '  -1.93%  '

$ws.Range('D51').Value = "'0.05957"
$ws.Range('E51').Value = '  +0.90%  '

